$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.082103155963927699
$ws.Range("B1").Value = -0.082103156380908524
$ws.Range("A2").Value = -0.020392998586153895
$ws.Range("B2").Value = 0.02039299817222703
$ws.Range("A3").Value = -0.049602940957051325
$ws.Range("B3").Value = 0.049602940529234647
$ws.Range("A4").Value = -0.040120672343372195
$ws.Range("B4").Value = 0.040120671903458978
